$wb = $excel.ActiveWorkbook

# --- Rename sheets (new timestamped names) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16512556056162512"
$wb.Worksheets.Item(2).Name = "NB_TO-16512556077562475"
$wb.Worksheets.Item(3).Name = "RS_TO-1651255607758249"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512556078212523"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512556078982496"

# --- Sheet 1 (GNG) stimulus file names ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512556055762498.csv"
$ws1.Range("B3").Value = "GNG_stims-16512556055992508.csv"
$ws1.Range("B4").Value = "go_stims-16512556056002512.csv"
$ws1.Range("B5").Value = "GNG_stims-16512556056152499.csv"

# --- Sheet 2 (NB) stimulus file names ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16512556075662482.csv"
$ws2.Range("B3").Value = "TB-1651255607739249.csv"
$ws2.Range("B4").Value = "OB-16512556060712495.csv"
$ws2.Range("B5").Value = "OB-16512556060462477.csv"
$ws2.Range("B6").Value = "TB-1651255606980251.csv"
$ws2.Range("B7").Value = "ZB-match_0-1651255605871251.csv"
$ws2.Range("B8").Value = "ZB-match_1-16512556056312494.csv"
$ws2.Range("B9").Value = "OB-1651255606414249.csv"
$ws2.Range("B10").Value = "ZB-match_1-16512556060142522.csv"

# --- Sheet 4 (TOL) stimulus file names ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512556077872493.csv"
$ws4.Range("B3").Value = "ZM_stims-16512556077652495.csv"
$ws4.Range("B4").Value = "MM_stims-165125560780325.csv"
$ws4.Range("B5").Value = "ZM_stims-16512556077882524.csv"
$ws4.Range("B6").Value = "MM_stims-1651255607819247.csv"
$ws4.Range("B7").Value = "ZM_stims-16512556078042517.csv"

# --- Sheet 5 (vSAT) stimulus file names ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1651255607883254.csv"
$ws5.Range("B3").Value = "SAT_stims-16512556078512523.csv"
$ws5.Range("B4").Value = "SAT_stims-16512556078272586.csv"
$ws5.Range("B5").Value = "vSAT_stims-16512556078662486.csv"
